$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift the existing ExpPoints column (C) to the right by 5 columns (to H),
#     inserting 5 new blank (but styled, where applicable) columns C:G.
#     This both relocates the old ExpPoints header/values to H and makes the
#     new header cells C1:G1 inherit the bold/centered/bordered style that
#     C1 ("ExpPoints") used to have. -4161 == xlShiftToRight
$ws.Range("C1:G21").Insert(-4161)

# --- Header row: label the five new columns; H1 already holds "ExpPoints" ---
$ws.Range("C1").Value = "WIN"
$ws.Range("D1").Value = "TOP4"
$ws.Range("E1").Value = "TOP5"
$ws.Range("F1").Value = "TOP6"
$ws.Range("G1").Value = "RELEGATION"

# --- Data rows: new team ranking + ExpPoints values; extra columns left blank ---
$teams = @(
    "Barcelona",
    "Real Madrid",
    "Atlético de Madrid",
    "Villarreal",
    "Real Betis",
    "Athletic Club",
    "Rayo Vallecano",
    "Celta de Vigo",
    "Espanyol",
    "Osasuna",
    "Real Sociedad",
    "Getafe",
    "Alavés",
    "Valencia",
    "Sevilla",
    "Mallorca",
    "Elche",
    "Levante",
    "Girona",
    "Real Oviedo"
)

$expPoints = @(
    83.33643001848033,
    82.73971312075047,
    72.27103819360818,
    63.77873058657867,
    57.86024465673832,
    54.47863749405111,
    50.92866782588134,
    49.34186892272135,
    46.88678366259734,
    46.49735567632037,
    46.14289385218007,
    45.69542048988172,
    45.21080254354546,
    44.29408897595547,
    43.86865846303058,
    42.58120123678579,
    41.52065142159509,
    38.28166269662648,
    35.22910093644742,
    34.19249272408548
)

for ($i = 0; $i -lt $teams.Length; $i++) {
    $row = $i + 2

    $ws.Cells.Item($row, 2).Value = $teams[$i]
    $ws.Cells.Item($row, 8).Value = $expPoints[$i]

    # Columns C..G (3..7) stay empty placeholders, same as the source sheet
    # (stored as an empty text value rather than being omitted entirely).
    for ($col = 3; $col -le 7; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        $cell.Value = "'"
        $cell.Style = $ws.Cells.Item($row, 1).Style
    }
}
